# Add data for 2022-08-11 (extends the "through August 02" running month
# column to "through August 03") to the carjacking-by-neighborhood-by-month
# workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet and update the running-total header label.
$ws.Name = "Through 2022-08-03"
$ws.Range("B1").Value = "August 2022 (through August 03)"

# Cell-level data updates (row => neighborhood).
# Austin (row 2)
$ws.Range("B2").Value = 2
$ws.Range("R2").Value = 1
$ws.Range("BF2").Value = 1

# Grand Crossing (row 3)
$ws.Range("B3").Value = 1
$ws.Range("D3").Value = 5

# Garfield Park (row 5)
$ws.Range("AH5").Value = 1

# Englewood (row 6)
$ws.Range("B6").Value = 1

# Chatham (row 9)
$ws.Range("J9").Value = 2
$ws.Range("R9").Value = 2

# West Town (row 12)
$ws.Range("B12").Value = 1
$ws.Range("AP12").Value = 1

# Roseland (row 13)
$ws.Range("B13").Value = 1

# Auburn Gresham (row 15)
$ws.Range("AH15").Value = 1

# Douglas (row 17)
$ws.Range("R17").Value = 1

# Little Italy, UIC (row 29)
$ws.Range("B29").Value = 4

# East Village (row 30)
$ws.Range("R30").Value = 1

# Calumet Heights (row 32)
$ws.Range("B32").Value = 1

# Washington Park (row 34)
$ws.Range("AP34").Value = 1

# West Pullman (row 37)
$ws.Range("B37").Value = 1

# South Chicago (row 39)
$ws.Range("R39").Value = 1
$ws.Range("AP39").Value = 3

# Armour Square (row 59)
$ws.Range("R59").Value = 1

# Avondale (row 61)
$ws.Range("J61").Value = 1

# Brighton Park (row 63)
$ws.Range("AH63").Value = 1

# Bucktown (row 64)
$ws.Range("J64").Value = 1

# Greektown (row 72)
$ws.Range("R72").Value = 1

# Hyde Park (row 74)
$ws.Range("J74").Value = 1

$wb.Save()
